# Apply the NPC configuration change:
# - Set column C (Public) to TRUE for rows 15 through 32
# - Update the active selection to C15:C38 with active cell C15

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 15; $row -le 32; $row++) {
    $ws.Cells.Item($row, 3).Value = $true
}

$ws.Range("C15:C38").Select()
